$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Comentarios": two cells in the fecha_comentario (I) column were
# mistakenly left with the date-only format (s=3); re-apply the date-time
# format (s=2) that the rest of that column uses.
# ---------------------------------------------------------------------------
$wsComentarios = $wb.Worksheets.Item("Comentarios")
$wsComentarios.Cells.Item(438, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsComentarios.Cells.Item(516, 9).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------------
# Sheet "Comentarios": a brand-new TikTok comment was scraped and inserted
# as row 638, pushing every following row down by one.
# ---------------------------------------------------------------------------
$wsComentarios.Rows.Item(638).Insert()

$wsComentarios.Cells.Item(638, 1).Value = 2
$wsComentarios.Cells.Item(638, 2).Value = "TikTok"
$wsComentarios.Cells.Item(638, 3).Value = "https://vt.tiktok.com/ZSfcQWN3t/"
$wsComentarios.Cells.Item(638, 4).Value = "https://vt.tiktok.com/ZSfcQWN3t/"
$wsComentarios.Cells.Item(638, 5).Value = ""
$wsComentarios.Cells.Item(638, 6).Value = "los quiero todos!"
$wsComentarios.Cells.Item(638, 7).Value = 1763678355
$wsComentarios.Cells.Item(638, 8).Value = 45981.94392361111
$wsComentarios.Cells.Item(638, 8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsComentarios.Cells.Item(638, 9).Value = 45981
$wsComentarios.Cells.Item(638, 9).NumberFormat = "YYYY-MM-DD"
$wsComentarios.Cells.Item(638, 10).Value = "22:39:15"
$wsComentarios.Cells.Item(638, 11).Value = 1
$wsComentarios.Cells.Item(638, 12).Value = 0
$wsComentarios.Cells.Item(638, 13).Value = $false
$wsComentarios.Cells.Item(638, 14).Value = "https://www.tiktok.com/@"
$wsComentarios.Cells.Item(638, 15).Value = ""
$wsComentarios.Cells.Item(638, 16).Value = "{'videoWebUrl': 'https://www.tiktok.com/@alpinacol/video/7574859079013633287', 'submittedVideoUrl': 'https://vt.tiktok.com/ZSfcQWN3t/', 'input': 'https://vt.tiktok.com/ZSfcQWN3t/', 'cid': '7574940790661808914', 'createTime': 1763678355, 'createTimeISO': '2025-11-20T22:39:15.000Z', 'text': 'los quiero todos!', 'diggCount': 1, 'likedByAuthor': False, 'pinnedByAuthor': False, 'repliesToId': None, 'replyCommentTotal': 0, 'uid': '6797815375695430662', 'uniqueId': 'marishoes', 'avatarThumbnail': 'http"

# ---------------------------------------------------------------------------
# Sheet "Resumen_Posts": refreshed comment/like totals for the TikTok post
# ---------------------------------------------------------------------------
$wsResumen = $wb.Worksheets.Item("Resumen_Posts")
$wsResumen.Cells.Item(3, 4).Value = 395
$wsResumen.Cells.Item(3, 5).Value = 322

# ---------------------------------------------------------------------------
# Sheet "Stats_Plataforma": refreshed platform-wide comment/like totals
# ---------------------------------------------------------------------------
$wsStats = $wb.Worksheets.Item("Stats_Plataforma")
$wsStats.Cells.Item(2, 3).Value = 688
$wsStats.Cells.Item(2, 5).Value = 1136
